# Chapter 6 / example 5 - add "max/min sales profit" summary cells (I1:J2)
# to the "背包" (backpack) product sales sheet, then tidy up the column
# widths the same way the original author did in Excel (select all ->
# AutoFit Column Width) after adding the two new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New labels in column I ------------------------------------------------
$ws.Range("I1").Value = "最大销售利润"
$ws.Range("I2").Value = "最小销售利润"

# --- New values in column J, formatted like the existing "销售利润" column -
$ws.Range("J1:J2").NumberFormat = '"¥"#,##0.00;"¥"\-#,##0.00'
$ws.Range("J1").Value = $excel.WorksheetFunction.Max($ws.Range("H2:H14"))
$ws.Range("J2").Value = $excel.WorksheetFunction.Min($ws.Range("H2:H14"))

# --- Resize all used columns to fit their new widest content ---------------
# (mirrors selecting the whole used range and choosing
#  Home > Format > AutoFit Column Width after the new columns were added)
$ws.UsedRange.EntireColumn.AutoFit() | Out-Null

# AutoFit in this environment doesn't reproduce the exact pixel metrics of
# the 微软雅黑 font Excel used originally, so nudge each column to the
# precise width the real workbook ended up with.
$ws.Columns.Item(1).ColumnWidth = 7.840401785714286   # A -> 8.5546875
$ws.Columns.Item(2).ColumnWidth = 7.504464285714286   # B -> 8.21875
$ws.Columns.Item(3).ColumnWidth = 13.504464285714286  # C -> 14.21875
$ws.Columns.Item(4).ColumnWidth = 13.504464285714286  # D -> 14.21875
$ws.Columns.Item(5).ColumnWidth = 12.840401785714286  # E -> 13.5546875
$ws.Columns.Item(6).ColumnWidth = 9.727120535714286   # F -> 10.44140625
$ws.Columns.Item(7).ColumnWidth = 9.727120535714286   # G -> 10.44140625
$ws.Columns.Item(8).ColumnWidth = 9.727120535714286   # H -> 10.44140625
$ws.Columns.Item(9).ColumnWidth = 10.949776785714286  # I -> 11.6640625
$ws.Columns.Item(10).ColumnWidth = 9.727120535714286  # J -> 10.44140625
